$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Todo ")
$ws2 = $wb.Worksheets.Item("Defect")

# --- "Todo " sheet -------------------------------------------------
# Row 4 ("Enemy attack") is no longer being worked on right now -
# clear out its PIC / Status.
$ws1.Range("E4").ClearContents()
$ws1.Range("H4").ClearContents()

# Row 6 ("UI" / village screen) is finished - assign PIC + Status.
$ws1.Range("E6").Value = "Fish"
$ws1.Range("H6").Value = "In progress"

# Row 9 ("level system") PIC / Status cleared too.
$ws1.Range("E9").ClearContents()
$ws1.Range("H9").ClearContents()

# --- "Defect" sheet -------------------------------------------------
# Remove the last two defect rows (4 & 5) entirely - they were
# resolved/duplicated into the Todo sheet above.
$ws2.Rows.Item(5).Delete()
$ws2.Rows.Item(4).Delete()

# Column A no longer needs to fit the long "Unusual Attack" /
# "Player pass through terrain" text, so shrink it back down.
$ws2.Columns.Item(1).ColumnWidth = 18.449776785714285

# --- selections -------------------------------------------------
$ws2.Select()
$ws2.Range("C24").Select()

$ws1.Select()
$ws1.Range("H6").Select()
